# Add 4 new rows (154-157) to Sheet1, continuing the existing daily series.
# Each new row replicates the values/formatting of the last existing row (153),
# except for column A, which holds the next sequential date serial number.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$sourceRow = 153
$firstNewRow = 154
$lastNewRow = 157

for ($r = $firstNewRow; $r -le $lastNewRow; $r++) {
    # Copy the whole source row (values + formatting) into the new row.
    $ws.Range("A$sourceRow`:J$sourceRow").Copy($ws.Range("A$r`:J$r"))

    # Column A is the date series; bump it by one day per new row.
    $ws.Cells.Item($r, 1).Value2 = 45709 + ($r - $sourceRow)
}

$wb.Save()
